$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column holds plain text (e.g. "1.001"); Excel would normally
# auto-detect such literals as numbers, so pre-format the cells we are about
# to rewrite as Text to keep them as strings, matching the source data.
$textRows = @(4,5,6,7,8,9,10,11,14,15,16,18,19,20,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,44,45,47,48,49,50)
foreach ($r in $textRows) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

# Updated Price values
$ws.Range("D2").Value = "30.166.55"
$ws.Range("D3").Value = "1.849.45"
$ws.Range("D4").Value = "1.001"
$ws.Range("D5").Value = "235.18"
$ws.Range("D6").Value = "1.001"
$ws.Range("D7").Value = "0.4702"
$ws.Range("D8").Value = "0.2885"
$ws.Range("D9").Value = "0.06523"
$ws.Range("D10").Value = "21.73"
$ws.Range("D11").Value = "0.07949"
$ws.Range("D13").Value = "1.856.33"
$ws.Range("D14").Value = "5.078"
$ws.Range("D15").Value = "0.6735"
$ws.Range("D16").Value = "265.76"
$ws.Range("D17").Value = "30.139.58"
$ws.Range("D18").Value = "13.58"
$ws.Range("D19").Value = "1.001"
$ws.Range("D20").Value = "0.000007545"
$ws.Range("D21").Value = "2.101.70"
$ws.Range("D23").Value = "5.199"
$ws.Range("D24").Value = "6.127"
$ws.Range("D25").Value = "166.65"
$ws.Range("D26").Value = "9.143"
$ws.Range("D27").Value = "18.77"
$ws.Range("D28").Value = "1.923"
$ws.Range("D29").Value = "1.393"
$ws.Range("D30").Value = "0.09824"
$ws.Range("D31").Value = "1.466"
$ws.Range("D32").Value = "4.259"
$ws.Range("D33").Value = "3.985"
$ws.Range("D34").Value = "0.04680"
$ws.Range("D35").Value = "1.114"
$ws.Range("D36").Value = "0.6961"
$ws.Range("D37").Value = "2.709"
$ws.Range("D38").Value = "0.01857"
$ws.Range("D39").Value = "2.600"
$ws.Range("D40").Value = "6.326"
$ws.Range("D41").Value = "73.09"
$ws.Range("D42").Value = "1.926"
$ws.Range("D43").Value = "1.001"
$ws.Range("D44").Value = "0.8362"
$ws.Range("D45").Value = "103.03"
$ws.Range("D47").Value = "942.01"
$ws.Range("D48").Value = "9.156"
$ws.Range("D49").Value = "6.991"
$ws.Range("D50").Value = "33.75"

# Updated Volume(1h) values
$ws.Range("E2").Value = "  -0.33%  "
$ws.Range("E3").Value = "  -0.84%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("E5").Value = "  +0.14%  "
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("E7").Value = "  +0.55%  "
$ws.Range("E8").Value = "  +1.84%  "
$ws.Range("E9").Value = "  +0.34%  "
$ws.Range("E10").Value = "  +2.25%  "
$ws.Range("E11").Value = "  +1.21%  "
$ws.Range("E12").Value = "  +0.28%  "
$ws.Range("E13").Value = "  -0.51%  "
$ws.Range("E14").Value = "  -0.15%  "
$ws.Range("E15").Value = "  +0.23%  "
$ws.Range("E16").Value = "  -4.93%  "
$ws.Range("E17").Value = "  -0.39%  "
$ws.Range("E18").Value = "  +7.21%  "
$ws.Range("E19").Value = "  +0.19%  "
$ws.Range("E20").Value = "  +3.95%  "
$ws.Range("E21").Value = "  -0.48%  "
$ws.Range("E22").Value = "  +0.17%  "
$ws.Range("E23").Value = "  -4.99%  "
$ws.Range("E24").Value = "  -0.18%  "
$ws.Range("E25").Value = "  +1.10%  "
$ws.Range("E26").Value = "  -0.36%  "
$ws.Range("E27").Value = "  -1.51%  "
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("E29").Value = "  +1.32%  "
$ws.Range("E30").Value = "  +2.06%  "
$ws.Range("E31").Value = "  -0.66%  "
$ws.Range("E32").Value = "  -2.77%  "
$ws.Range("E33").Value = "  -2.55%  "
$ws.Range("E34").Value = "  -0.22%  "
$ws.Range("E35").Value = "  -0.24%  "
$ws.Range("E36").Value = "  -1.27%  "
$ws.Range("E37").Value = "  -0.69%  "
$ws.Range("E38").Value = "  +0.44%  "
$ws.Range("E39").Value = "  +2.59%  "
$ws.Range("E40").Value = "  +1.33%  "
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("E42").Value = "  -0.70%  "
$ws.Range("E43").Value = "  +0.07%  "
$ws.Range("E44").Value = "  -1.01%  "
$ws.Range("E45").Value = "  -0.58%  "
$ws.Range("E46").Value = "  -1.15%  "
$ws.Range("E47").Value = "  +0.67%  "
$ws.Range("E48").Value = "  -0.04%  "
$ws.Range("E49").Value = "  -2.38%  "
$ws.Range("E50").Value = "  -0.72%  "
